$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (bold)
$j12 = $ws.Range("J12")
$j12.Formula = "=AVERAGE(J2:J11)"
$j12.Font.Bold = $true

# Row 14-17: summary stats
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Format B14 (bold, size 12, vertically centered), then replicate to B15:B17
$c = $ws.Range("B14")
$c.Font.Bold = $true
$c.Font.Size = 12
$c.VerticalAlignment = -4108

$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)

$ws.Range("A14:B17").RowHeight = 15.6

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("A14:B17").Select()
